$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 4.5
$ws.Range("H2").Value = 4.1
$ws.Range("I2").Value = 1.67
$ws.Range("K2").Value = 2.4
$ws.Range("L2").Value = 2.2
$ws.Range("O2").Value = 1.14
$ws.Range("P2").Value = 5.5
$ws.Range("Q2").Value = 1.53
$ws.Range("R2").Value = 2.4
$ws.Range("S2").Value = 1.29
$ws.Range("T2").Value = 3.5
$ws.Range("W2").Value = 17
$ws.Range("X2").Value = 26
$ws.Range("AD2").Value = 8
$ws.Range("AH2").Value = 9.5
$ws.Range("AM2").Value = 21
$ws.Range("AN2").Value = 6.5
$ws.Range("AS2").Value = 151
$ws.Range("AT2").Value = 3.5
$ws.Range("AY2").Value = 8.5
$ws.Range("AZ2").Value = 17
$ws.Range("BC2").Value = 101
$ws.Range("G3").Value = 1.2
$ws.Range("H3").Value = 6.5
$ws.Range("Y3").Value = 10
$ws.Range("AB3").Value = 29
$ws.Range("AG3").Value = 351
$ws.Range("AJ3").Value = 29
$ws.Range("AQ3").Value = 12
$ws.Range("AY3").Value = 51
$ws.Range("BA3").Value = 251
$ws.Range("BC3").Value = 351
$ws.Range("G4").Value = 2.1
$ws.Range("H4").Value = 3.6
$ws.Range("I4").Value = 3.25
$ws.Range("J4").Value = 2.75
$ws.Range("L4").Value = 3.75
$ws.Range("Q4").Value = 1.8
$ws.Range("R4").Value = 2
$ws.Range("S4").Value = 1.36
$ws.Range("T4").Value = 3
$ws.Range("AA4").Value = 17
$ws.Range("AC4").Value = 12
$ws.Range("AH4").Value = 11
$ws.Range("AI4").Value = 17
$ws.Range("AK4").Value = 34
$ws.Range("AL4").Value = 23
$ws.Range("AP4").Value = 21
$ws.Range("AQ4").Value = 41
$ws.Range("AT4").Value = 3
$ws.Range("AY4").Value = 17
$ws.Range("G5").Value = 1.7
$ws.Range("I5").Value = 5.25
$ws.Range("J5").Value = 2.4
$ws.Range("K5").Value = 1.95
$ws.Range("U5").Value = 2.38
$ws.Range("V5").Value = 1.53
$ws.Range("AC5").Value = 6.5
$ws.Range("AM5").Value = 51
$ws.Range("AO5").Value = 9.5
$ws.Range("BD5").Value = 126
$ws.Range("I6").Value = 3.6
$ws.Range("M6").Value = 1.06
$ws.Range("N6").Value = 10
$ws.Range("Q6").Value = 1.93
$ws.Range("R6").Value = 1.93
$ws.Range("S6").Value = 1.4
$ws.Range("T6").Value = 2.75
$ws.Range("U6").Value = 1.73
$ws.Range("V6").Value = 2
$ws.Range("W6").Value = 8
$ws.Range("AB6").Value = 26
$ws.Range("AH6").Value = 11
$ws.Range("AJ6").Value = 13
$ws.Range("AL6").Value = 29
$ws.Range("AT6").Value = 2.75
$ws.Range("AU6").Value = 8
$ws.Range("BA6").Value = 67
$ws.Range("BC6").Value = 201
$ws.Range("Q7").Value = 1.8
$ws.Range("R7").Value = 2
$ws.Range("G8").Value = 3.25
$ws.Range("Y8").Value = 12
$ws.Range("AH8").Value = 8
$ws.Range("G9").Value = 2.25
$ws.Range("H9").Value = 3.3
$ws.Range("I9").Value = 3.1
$ws.Range("J9").Value = 3.1
$ws.Range("K9").Value = 2
$ws.Range("S9").Value = 1.5
$ws.Range("T9").Value = 2.5
$ws.Range("U9").Value = 1.91
$ws.Range("V9").Value = 1.8
$ws.Range("AA9").Value = 21
$ws.Range("AC9").Value = 8.5
$ws.Range("AG9").Value = 401
$ws.Range("AH9").Value = 8
$ws.Range("AP9").Value = 26
$ws.Range("AT9").Value = 2.5
$ws.Range("AV9").Value = 67
$ws.Range("BB9").Value = 101
$ws.Range("G10").Value = 1.85
$ws.Range("I10").Value = 4.2
$ws.Range("N10").Value = 9
$ws.Range("AH10").Value = 11
$ws.Range("AO10").Value = 10
$ws.Range("AQ10").Value = 34
$ws.Range("AS10").Value = 151
$ws.Range("G11").Value = 1.75
$ws.Range("H11").Value = 3.6
$ws.Range("I11").Value = 4.5
$ws.Range("J11").Value = 2.38
$ws.Range("L11").Value = 4.75
$ws.Range("Q11").Value = 1.93
$ws.Range("R11").Value = 1.93
$ws.Range("W11").Value = 7
$ws.Range("X11").Value = 8.5
$ws.Range("Z11").Value = 13
$ws.Range("AC11").Value = 11
$ws.Range("AD11").Value = 7
$ws.Range("AI11").Value = 23
$ws.Range("AJ11").Value = 15
$ws.Range("AK11").Value = 51
$ws.Range("AO11").Value = 9
$ws.Range("AQ11").Value = 29
$ws.Range("AX11").Value = 6.5
$ws.Range("AY11").Value = 23
$ws.Range("G13").Value = 2.15
$ws.Range("H13").Value = 3.05
$ws.Range("J13").Value = 2.67
$ws.Range("K13").Value = 2.07
$ws.Range("L13").Value = 3.8
$ws.Range("O13").Value = 1.33
$ws.Range("P13").Value = 2.8
$ws.Range("U13").Value = 1.75
$ws.Range("V13").Value = 1.87
$ws.Range("W13").Value = 7.2
$ws.Range("X13").Value = 10.25
$ws.Range("Z13").Value = 21
$ws.Range("AA13").Value = 18
$ws.Range("AB13").Value = 29
$ws.Range("AD13").Value = 6
$ws.Range("AE13").Value = 14
$ws.Range("AF13").Value = 65
$ws.Range("AI13").Value = 17.5
$ws.Range("AJ13").Value = 11.5
$ws.Range("AK13").Value = 50
$ws.Range("AN13").Value = 4.05
$ws.Range("AO13").Value = 10.75
$ws.Range("AP13").Value = 17.5
$ws.Range("AQ13").Value = 40
$ws.Range("AR13").Value = 65
$ws.Range("AS13").Value = 200
$ws.Range("AT13").Value = 2.57
$ws.Range("AU13").Value = 6.6
$ws.Range("AV13").Value = 55
$ws.Range("AX13").Value = 5.3
$ws.Range("AY13").Value = 18.5
$ws.Range("AZ13").Value = 24
$ws.Range("BA13").Value = 90
$ws.Range("BB13").Value = 120
$ws.Range("BC13").Value = 300
$ws.Range("G14").Value = 2.47
$ws.Range("I14").Value = 2.62
$ws.Range("J14").Value = 3.05
$ws.Range("L14").Value = 3.15
$ws.Range("T14").Value = 2.95
$ws.Range("W14").Value = 10.25
$ws.Range("X14").Value = 14.5
$ws.Range("Z14").Value = 29
$ws.Range("AA14").Value = 19
$ws.Range("AB14").Value = 23
$ws.Range("AH14").Value = 10.75
$ws.Range("AI14").Value = 15.5
$ws.Range("AL14").Value = 20
$ws.Range("AM14").Value = 23
$ws.Range("AP14").Value = 18
$ws.Range("AQ14").Value = 55
$ws.Range("AT14").Value = 2.95
$ws.Range("AX14").Value = 4.85
$ws.Range("AY14").Value = 14
$ws.Range("AZ14").Value = 18.5
$ws.Range("BB14").Value = 75
$ws.Range("BC14").Value = 175
